$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the emoji, replace with a plain-text smiley.
# ------------------------------------------------------------------
$find = $d.Content.Find
$find.Execute("🤓", $true, $false, $false, $false, $false, $true, 1, $false, ":)", 2)

# ------------------------------------------------------------------
# 2. Add a new bullet under "Design / Development / Security",
#    right before the "As a BAE Systems principal engineer..." bullet.
# ------------------------------------------------------------------
$count = $d.Paragraphs.Count
$targetIndex = 0
for ($i = 1; $i -le $count; $i++) {
    $txt = $d.Paragraphs.Item($i).Range.Text
    if ($txt -like "As a BAE Systems principal engineer, consulted on multiple projects*") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -gt 0) {
    $anchorPara = $d.Paragraphs.Item($targetIndex)
    # Inserting a paragraph *before* this bullet copies its paragraph
    # style (Compact) and numbering (numId 1007) onto the new paragraph.
    $anchorPara.Range.InsertParagraphBefore()

    $newPara = $d.Paragraphs.Item($targetIndex)
    $rng = $newPara.Range
    $rng.Collapse(1)
    $startPos = $rng.Start

    $part1 = "At USSTRATCOM, designed a visual grammar in order to create detailed process documentation. The documentation is intended to help"
    $part2 = " "
    $part3 = "architect sustainable solutions"
    $part4 = " "
    $part5 = "in an effort to improve efficiency and mission readiness."

    $rng.InsertAfter($part1 + $part2 + $part3 + $part4 + $part5)

    $boldStart = $startPos + $part1.Length + $part2.Length
    $boldEnd = $boldStart + $part3.Length
    $boldRng = $d.Range($boldStart, $boldEnd)
    $boldRng.Bold = 1
}
